# Generate Report for handback
# This script updates the localization-status workbook to reflect that the
# zh-cn and de-de handoffs have now been handed back: it records the target
# (translated) file + handback file for the a.md.md row, updates the
# "Latest Handback DateTime" for both files, and flips the overall status
# text from "Ready for handoff" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1. Update the status text everywhere it currently appears.
# ---------------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")
if ($ovw.Range("B2").Value() -eq $oldStatus) { $ovw.Range("B2").Value = $newStatus }
if ($ovw.Range("C2").Value() -eq $oldStatus) { $ovw.Range("C2").Value = $newStatus }
if ($ovw.Range("B3").Value() -eq $oldStatus) { $ovw.Range("B3").Value = $newStatus }
if ($ovw.Range("C3").Value() -eq $oldStatus) { $ovw.Range("C3").Value = $newStatus }

$zhcn = $wb.Worksheets.Item("zh-cn")
if ($zhcn.Range("B2").Value() -eq $oldStatus) { $zhcn.Range("B2").Value = $newStatus }
if ($zhcn.Range("B3").Value() -eq $oldStatus) { $zhcn.Range("B3").Value = $newStatus }

$dede = $wb.Worksheets.Item("de-de")
if ($dede.Range("B2").Value() -eq $oldStatus) { $dede.Range("B2").Value = $newStatus }
if ($dede.Range("B3").Value() -eq $oldStatus) { $dede.Range("B3").Value = $newStatus }

# ---------------------------------------------------------------------
# 2. zh-cn sheet: fill in "Latest Target File" (E) / "Latest Handback
#    File" (F) hyperlinks for the two tracked markdown files, and stamp
#    the "Latest Handback DateTime" (G) now that the handback happened.
# ---------------------------------------------------------------------
$zhSourceUrl = "https://github.com/OpenLocalizationTest/oltest/blob/8763f0fbce2541b697cf075d5ece242c51ec3288/e2e/a.md.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6812aef6f9aea804fb060b0052d33efd177a09a0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf"
$zhXlfName = "a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf"

$zhcn.Hyperlinks.Add($zhcn.Range("E2"), $zhSourceUrl, "", "", "a.md.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), $zhXlfUrl, "", "", $zhXlfName) | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("E3"), $zhSourceUrl, "", "", "a.md.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), $zhXlfUrl, "", "", $zhXlfName) | Out-Null

$zhcn.Range("G2").Value = "2016-01-18 10:58:59"
$zhcn.Range("G3").Value = "2016-01-18 10:58:59"

# ---------------------------------------------------------------------
# 3. de-de sheet: same treatment.
# ---------------------------------------------------------------------
$deSourceUrl = "https://github.com/OpenLocalizationTest/oltest/blob/8763f0fbce2541b697cf075d5ece242c51ec3288/e2e/a.md.md"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e6299fae389e705aa9e17aaf85bdd45ba59feca9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf"
$deXlfName = "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf"

$dede.Hyperlinks.Add($dede.Range("E2"), $deSourceUrl, "", "", "a.md.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("F2"), $deXlfUrl, "", "", $deXlfName) | Out-Null
$dede.Hyperlinks.Add($dede.Range("E3"), $deSourceUrl, "", "", "a.md.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("F3"), $deXlfUrl, "", "", $deXlfName) | Out-Null

$dede.Range("G2").Value = "2016-01-18 10:59:16"
$dede.Range("G3").Value = "2016-01-18 10:59:16"
